# Pump current (Ip) calculation sheet added; ADC reference sheet renamed.
$wb = $excel.ActiveWorkbook

# --- Rename the original sheet "Tabelle1" -> "ADC" ---------------------
$adc = $wb.Worksheets.Item(1)
$adc.Name = "ADC"

# --- Add the new "Ip" sheet right after "ADC" ---------------------------
$ip = $wb.Worksheets.Add($null, $adc)
$ip.Name = "Ip"

# Row 1 - Ua
$ip.Range("A1").Value = "Ua"
$ip.Range("B1").Value = 0.51
$ip.Range("C1").Formula = "=B1*1000"

# Row 2 - Ua_cal
$ip.Range("A2").Value = "Ua_cal"
$ip.Range("B2").Value = 1.5
$ip.Range("C2").Formula = "=B2*1000"

# Row 3 - Diff
$ip.Range("A3").Value = "Diff"
$ip.Range("B3").Formula = "=B1-B2"
$ip.Range("C3").Formula = "=C1-C2"

# Row 4 - amp
$ip.Range("A4").Value = "amp"
$ip.Range("B4").Value = 8
$ip.Range("C4").Value = 8

# Row 5 - Rshunt
$ip.Range("A5").Value = "Rshunt"
$ip.Range("B5").Value = 61.9
$ip.Range("C5").Value = 61.9

# Row 7 - Ip (pump current)
$ip.Range("A7").Value = "Ip"
$ip.Range("B7").Formula = "=((B1-B2)/(B4*B5))*1000"
$ip.Range("B7").NumberFormat = "0.000"
$ip.Range("C7").Formula = "=(C1-C2)/(C4*C5)"

# Selection / view state: "Ip" becomes the active sheet, cell B1 selected
[void]$ip.Range("B1").Select()

# Page margins on the new sheet are metric (2 cm) rather than the 0.75in default
$ip.PageSetup.TopMargin = 56.692913399999995
$ip.PageSetup.BottomMargin = 56.692913399999995
